$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 18650
$ws.Range("E3").Value = 11807
$ws.Range("E4").Value = 12839
$ws.Range("E5").Value = 8604
$ws.Range("E6").Value = 16563
$ws.Range("E7").Value = 3264
$ws.Range("E8").Value = 12426
$ws.Range("E9").Value = 1083
$ws.Range("E10").Value = 9518
$ws.Range("E11").Value = 18685
$ws.Range("E12").Value = 3876
$ws.Range("E13").Value = 7879
